# Update NATMI LR-pair TPM output (Epo-Epor) with newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (ECs -> Epo -> Epor -> ECs)
$ws.Range("E2").Value  = 1
$ws.Range("F2").Value  = 0.3333333333333333
$ws.Range("G2").Value  = 0.05285566666666667
$ws.Range("H2").Value  = 0.158567
$ws.Range("I2").Value  = 0.5476571215423245
$ws.Range("J2").Value  = 0.5476571215423245
$ws.Range("K2").Value  = 3
$ws.Range("L2").Value  = 1
$ws.Range("M2").Value  = 0.4065259999999999
$ws.Range("N2").Value  = 1.219578
$ws.Range("Q2").Value  = 0.02148720274733333
$ws.Range("R2").Value  = 0.193384824726
$ws.Range("S2").Value  = 0.5476571215423245
$ws.Range("T2").Value  = 0.5476571215423245

# Row 3 (FAPs -> Epo -> Epor -> ECs)
$ws.Range("I3").Value  = 0.3108894545429426
$ws.Range("J3").Value  = 0.3108894545429427
$ws.Range("K3").Value  = 3
$ws.Range("L3").Value  = 1
$ws.Range("M3").Value  = 0.4065259999999999
$ws.Range("N3").Value  = 1.219578
$ws.Range("Q3").Value  = 0.01219767712133333
$ws.Range("R3").Value  = 0.109779094092
$ws.Range("S3").Value  = 0.3108894545429426
$ws.Range("T3").Value  = 0.3108894545429427

# Row 4 (MuSCs -> Epo -> Epor -> ECs)
$ws.Range("I4").Value  = 0.1414534239147328
$ws.Range("J4").Value  = 0.1414534239147328
$ws.Range("K4").Value  = 3
$ws.Range("L4").Value  = 1
$ws.Range("M4").Value  = 0.4065259999999999
$ws.Range("N4").Value  = 1.219578
$ws.Range("Q4").Value  = 0.005549892951999999
$ws.Range("R4").Value  = 0.04994903656799999
$ws.Range("S4").Value  = 0.1414534239147328
$ws.Range("T4").Value  = 0.1414534239147328
